$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2462006079027356
$ws.Range("C2").Value = 0.4589665653495441
$ws.Range("J2").Value = 0.02735562310030395
$ws.Range("O2").Value = 0.00303951367781155
$ws.Range("P2").Value = 0.1793313069908815
$ws.Range("S2").Value = 0.0851063829787234
$ws.Range("B3").Value = 0.03125
$ws.Range("C3").Value = 0.0375
$ws.Range("J3").Value = 0.0375
$ws.Range("P3").Value = 0.75625
$ws.Range("S3").Value = 0.1375
$ws.Range("J4").Value = 0.02127659574468085
$ws.Range("P4").Value = 0.7872340425531915
$ws.Range("S4").Value = 0.1914893617021277
$ws.Range("P5").Value = 0.8
$ws.Range("S5").Value = 0.2
$ws.Range("B6").Value = 0.05676855895196507
$ws.Range("D6").Value = 0.008733624454148471
$ws.Range("F6").Value = 0.08296943231441048
$ws.Range("J6").Value = 0.2445414847161572
$ws.Range("O6").Value = 0.03930131004366812
$ws.Range("Q6").Value = 0.1179039301310044
$ws.Range("R6").Value = 0.09606986899563319
$ws.Range("S6").Value = 0.3537117903930131
$ws.Range("B7").Value = 0.06779661016949153
$ws.Range("D7").Value = 0.02259887005649718
$ws.Range("E7").Value = 0.005649717514124294
$ws.Range("F7").Value = 0.06779661016949153
$ws.Range("J7").Value = 0.1016949152542373
$ws.Range("O7").Value = 0.02824858757062147
$ws.Range("Q7").Value = 0.1751412429378531
$ws.Range("R7").Value = 0.07909604519774012
$ws.Range("S7").Value = 0.4519774011299435
$ws.Range("B8").Value = 0.09310986964618249
$ws.Range("D8").Value = 0.01675977653631285
$ws.Range("E8").Value = 0.00186219739292365
$ws.Range("F8").Value = 0.05772811918063315
$ws.Range("J8").Value = 0.1024208566108007
$ws.Range("O8").Value = 0.0223463687150838
$ws.Range("Q8").Value = 0.143389199255121
$ws.Range("R8").Value = 0.111731843575419
$ws.Range("S8").Value = 0.4506517690875233
$ws.Range("B9").Value = 0.1551724137931035
$ws.Range("D9").Value = 0.01724137931034483
$ws.Range("E9").Value = 0.005747126436781609
$ws.Range("F9").Value = 0.08045977011494253
$ws.Range("J9").Value = 0.06896551724137931
$ws.Range("O9").Value = 0.005747126436781609
$ws.Range("Q9").Value = 0.132183908045977
$ws.Range("R9").Value = 0.103448275862069
$ws.Range("S9").Value = 0.4310344827586207
$ws.Range("B10").Value = 0.1091350040420372
$ws.Range("D10").Value = 0.02506063055780113
$ws.Range("E10").Value = 0.001616814874696847
$ws.Range("F10").Value = 0.06063055780113177
$ws.Range("J10").Value = 0.110751818916734
$ws.Range("O10").Value = 0.01616814874696847
$ws.Range("Q10").Value = 0.2077607113985449
$ws.Range("R10").Value = 0.09539207760711399
$ws.Range("S10").Value = 0.3734842360549717
$ws.Range("G11").Value = 0.1354581673306773
$ws.Range("J11").Value = 0.08366533864541832
$ws.Range("K11").Value = 0.1752988047808765
$ws.Range("L11").Value = 0.5737051792828686
$ws.Range("S11").Value = 0.03187250996015936
$ws.Range("G12").Value = 0.756578947368421
$ws.Range("J12").Value = 0.1710526315789474
$ws.Range("L12").Value = 0.05263157894736842
$ws.Range("S12").Value = 0.01973684210526316
$ws.Range("G13").Value = 0.6538461538461539
$ws.Range("J13").Value = 0.2307692307692308
$ws.Range("S13").Value = 0.1153846153846154
$ws.Range("F15").Value = 0.02597402597402598
$ws.Range("H15").Value = 0.1341991341991342
$ws.Range("I15").Value = 0.04761904761904762
$ws.Range("J15").Value = 0.3506493506493507
$ws.Range("K15").Value = 0.06060606060606061
$ws.Range("M15").Value = 0.01298701298701299
$ws.Range("N15").Value = 0.004329004329004329
$ws.Range("O15").Value = 0.05627705627705628
$ws.Range("S15").Value = 0.3073593073593073
$ws.Range("F16").Value = 0.02325581395348837
$ws.Range("H16").Value = 0.2465116279069768
$ws.Range("I16").Value = 0.08372093023255814
$ws.Range("J16").Value = 0.2883720930232558
$ws.Range("K16").Value = 0.08837209302325581
$ws.Range("M16").Value = 0.0186046511627907
$ws.Range("N16").Value = 0.004651162790697674
$ws.Range("O16").Value = 0.06511627906976744
$ws.Range("S16").Value = 0.1813953488372093
$ws.Range("F17").Value = 0.01895734597156398
$ws.Range("H17").Value = 0.2132701421800948
$ws.Range("I17").Value = 0.06635071090047394
$ws.Range("J17").Value = 0.4075829383886256
$ws.Range("K17").Value = 0.08767772511848342
$ws.Range("M17").Value = 0.02369668246445497
$ws.Range("N17").Value = 0.004739336492890996
$ws.Range("O17").Value = 0.07819905213270142
$ws.Range("S17").Value = 0.0995260663507109
$ws.Range("F18").Value = 0.008583690987124463
$ws.Range("H18").Value = 0.2489270386266094
$ws.Range("I18").Value = 0.09442060085836911
$ws.Range("J18").Value = 0.4034334763948498
$ws.Range("K18").Value = 0.06008583690987124
$ws.Range("M18").Value = 0.01716738197424893
$ws.Range("O18").Value = 0.05150214592274678
$ws.Range("S18").Value = 0.1158798283261803
$ws.Range("F19").Value = 0.02086438152011923
$ws.Range("H19").Value = 0.2302533532041729
$ws.Range("I19").Value = 0.07078986587183309
$ws.Range("J19").Value = 0.3673621460506706
$ws.Range("K19").Value = 0.08718330849478391
$ws.Range("M19").Value = 0.02309985096870343
$ws.Range("O19").Value = 0.06706408345752608
$ws.Range("S19").Value = 0.1333830104321908
